$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 321, shifting existing rows 321:415 down to 322:416
$ws.Rows("321:321").Insert()

# Populate the newly inserted row 321 with the new weekly price record
$ws.Cells.Item(321, 1).Value = 8
$ws.Cells.Item(321, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(321, 3).Value = "Coquimbo"
$ws.Cells.Item(321, 4).Value = 44841
$ws.Cells.Item(321, 5).Value = 4
$ws.Cells.Item(321, 6).Value = 100114013
$ws.Cells.Item(321, 7).Value = "Zanahoria"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 600
$ws.Cells.Item(321, 11).Value = 12800
$ws.Cells.Item(321, 12).Value = 13000
$ws.Cells.Item(321, 13).Value = 12900
$ws.Cells.Item(321, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(321, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(321, 16).Value = 645
$ws.Cells.Item(321, 17).Value = 20
$ws.Cells.Item(321, 18).Value = "Hortaliza"
